# Applies the commit "fixed legends and tables":
#   - rename the sheet from "Normalized cDNA Reads and Annotation" to "Sheet1"
#   - bump the sheet zoom from 90% to 110%
#   - move the active selection from B9 to D22
#   - shrink the first ten column widths slightly (their character-width
#     values all decrease by a small, consistent amount)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sheet name ---------------------------------------------------------
$ws.Name = "Sheet1"

# --- zoom ----------------------------------------------------------------
$excel.ActiveWindow.Zoom = 110

# keep gridlines visible (matches the original sheetView's showGridLines)
$excel.ActiveWindow.DisplayGridlines = $true

# --- selection -------------------------------------------------------------
$ws.Range("D22").Select()

# --- column widths ---------------------------------------------------------
# Target stored (OOXML) widths from the diff; ColumnWidth (Excel "characters"
# units) maps to stored width via: stored = ColumnWidth + 5/6, so subtract
# that offset before assigning. Columns G and H share one merged <col> range
# in the sheet (both carry the same width), so they get the same value.
$offset = 5/6

# A, B, C, D, E, F, G, H, I, J
$targetWidths = @(
    8.10204081632653,
    40.0918367346939,
    38.6071428571429,
    39.8214285714286,
    29.1581632653061,
    119.872448979592,
    315.612244897959,
    315.612244897959,
    137.423469387755,
    30.2397959183673
)

for ($i = 0; $i -lt $targetWidths.Length; $i++) {
    $col = $i + 1
    $ws.Columns.Item($col).ColumnWidth = $targetWidths[$i] - $offset
}
